# Populate Sheet1 with the manufacturer input rows (rows 6-14).
# Row 6 already has its part-number (A6); fill in the Manufacturer / Cost
# columns for it, then add the remaining rows below it.

function Set-InputRow($Sheet, $RowIndex, $PartNumber, $Manufacturer, $CostPrice) {
    $Sheet.Cells.Item($RowIndex, 1).Value = $PartNumber
    $Sheet.Cells.Item($RowIndex, 2).Value = $Manufacturer
    $Sheet.Cells.Item($RowIndex, 3).Value = $CostPrice
}

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

$inputRows = @(
    @{ PartNumber = 22356568; Manufacturer = "Rnd Manufacturer"; CostPrice = 22691 },
    @{ PartNumber = 34256354; Manufacturer = "Rnd Manufacturer"; CostPrice = 22691 },
    @{ PartNumber = 44637355; Manufacturer = "Rnd Manufacturer"; CostPrice = 22691 },
    @{ PartNumber = 55261434; Manufacturer = "Rnd Manufacturer"; CostPrice = 22691 },
    @{ PartNumber = 66464788; Manufacturer = "Rnd Manufacturer"; CostPrice = 22691 },
    @{ PartNumber = 77423423; Manufacturer = "Rnd Manufacturer"; CostPrice = 22691 },
    @{ PartNumber = 88888856; Manufacturer = "Rnd Manufacturer"; CostPrice = 22691 },
    @{ PartNumber = 94757647; Manufacturer = "Rnd Manufacturer"; CostPrice = 22691 },
    @{ PartNumber = 10342423; Manufacturer = "Rnd Manufacturer"; CostPrice = 22691 }
)

$startRow = 6
for ($i = 0; $i -lt $inputRows.Count; $i++) {
    $row = $inputRows[$i]
    $rowIndex = $startRow + $i
    Set-InputRow $ws1 $rowIndex $row.PartNumber $row.Manufacturer $row.CostPrice
}

# Remove the leftover test values from Sheet2 (rows 6-7), shrinking the
# sheet's used range back down to A1:F5.
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Rows.Item(6).Delete()
$ws2.Rows.Item(6).Delete()
